$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Row 3 updates
$ws.Range("A3").Value = 56495835
$ws.Range("S3").Value = 10
$ws.Range("AO3").Value = "1 substratenheter # Rönn"

# Row 4 updates
$ws.Range("A4").Value = 56495834
$ws.Range("S4").Value = 10
$ws.Range("AO4").Value = "1 substratenheter # Rönn"
